$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, shifting existing rows 5-7 down to 6-8.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new weekly price entry.
$ws.Range("A5").Value = 10
$ws.Range("B5").Value = "Vega Modelo de Temuco"
$ws.Range("C5").Value = "La Araucanía"
$ws.Range("D5").Value = 44567
$ws.Range("D5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100101
$ws.Range("H5").Value = "Berries"
$ws.Range("I5").Value = 100101004
$ws.Range("J5").Value = "Frambuesa"
$ws.Range("K5").Value = "Sin especificar"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 80
$ws.Range("N5").Value = 2400
$ws.Range("O5").Value = 2400
$ws.Range("P5").Value = 2400
$ws.Range("Q5").Value = "`$/envase 1 kilo"
$ws.Range("R5").Value = "Región de La Araucanía"
$ws.Range("S5").Value = 2400
$ws.Range("T5").Value = 1
